$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Add lecture 5" — new participation/attendance column G ("Assignment 4"
# header already existed in G1); mark the 1s for lecture-5 attendance, bump
# a few assignment-count cells in column C, and flag three students who
# dropped the course (DROP, with strikethrough formatting on their rows).
# ---------------------------------------------------------------------------

# New column G (lecture 5) attendance marks
$ws.Range("G2").Value  = 1
$ws.Range("G6").Value  = 1
$ws.Range("G10").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("G17").Value = 1
$ws.Range("G19").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("G23").Value = 1
$ws.Range("G26").Value = 1
$ws.Range("G27").Value = 1
$ws.Range("G28").Value = 1
$ws.Range("G32").Value = 1
$ws.Range("G33").Value = 1

# Column C (assignment count) updates / additions
$ws.Range("C6").Value  = 3
$ws.Range("C13").Value = 2
$ws.Range("C17").Value = 2
$ws.Range("C26").Value = 3
$ws.Range("C28").Value = 1
$ws.Range("C29").Value = 1
$ws.Range("C33").Value = 2

# Misc newly-filled cells
$ws.Range("F13").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("D29").Value = 1

# ---------------------------------------------------------------------------
# Dropped students: Covarrubis (row 7), Hackett (row 11), Trejo (row 31).
# Strikethrough the row's existing cells; row 7 additionally gets a "DROP"
# marker in column C.
# ---------------------------------------------------------------------------

$ws.Range("C7").Value = "DROP"

$ws.Range("A7:B7").Font.Strikethrough = $true
$ws.Range("C7:D7").Font.Strikethrough = $true

$ws.Range("A11:B11").Font.Strikethrough = $true
$ws.Range("D11:F11").Font.Strikethrough = $true

$ws.Range("A31:B31").Font.Strikethrough = $true
$ws.Range("D31:E31").Font.Strikethrough = $true

# ---------------------------------------------------------------------------
# Selection moved to E11
# ---------------------------------------------------------------------------
[void]$ws.Range("E11").Select()
